# The "HS" (Highest Score) column (G) stores not-out innings as text
# like "89*", "100*", etc. This pass strips the trailing "*" and rewrites
# those cells as plain numbers (89, 100, ...). Once nothing references the
# old "NN*" text any more, saving naturally drops those now-unused shared
# string entries, which shifts every later shared-string index down to
# fill the gaps left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = HS
    $val = $cell.Value()
    if ($val -ne $null -and $val -match '^\d+\*$') {
        $numPart = $val.Substring(0, $val.Length - 1)
        $num = $numPart -as [int]
        $cell.Value = $num
    }
}

# Matches the trailing view-state change in the diff: the active cell
# moves to G10 (and the window scrolls so row 75 is visible, though the
# scroll position itself is a cosmetic window property this host does not
# persist to the saved file).
[void]$ws.Range("G10").Select()
